$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the game-date label, which was previously stored as
# "2-27-2011-12" (a malformed mash-up of the folder name / old date
# string). NBA.com's boxscore date and the intended ISO date is
# 2012-02-27, so rewrite every data row (2 through 31) in column BF.
#
# Pre-format the cells as Text so Excel's auto-detection doesn't turn
# the ISO-looking string back into a date serial number - we want the
# literal text "2012-02-27" stored in the cell.
$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "2012-02-27"
}
